# Update the "Rules" worksheet: cell E8 previously held "Good Morning",
# change it to "GIT UPDATE" (this also makes "Good Morning" disappear from
# the shared-strings table since it becomes unused, while "GIT UPDATE" is
# appended as a new shared string).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = "GIT UPDATE"

# Reflect that E8 is the active/selected cell in the saved view.
$ws.Range("E8").Select()
